$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared strings: Sule -> a, Sule Myodaw Hall -> b, Yoke Shin Yone -> c
$ws.Range("A3").Value = "a"
$ws.Range("A5").Value = "b"
$ws.Range("A7").Value = "c"

# Update the selection to C3
$ws.Range("C3").Select()
